# Update the category-by-subject AE table: total accrual N changed from 2 to 3,
# so the header labels and every derived percentage must be refreshed.
#
# Plain "$range.Value = '66.67'" gets auto-coerced to a number by Excel's
# type inference (and pre-setting NumberFormat to force text leaves a stray
# style behind). Instead, write the text via a formula that evaluates to a
# string literal, then Copy + PasteSpecial(xlPasteValues) over itself so the
# formula collapses to a plain shared-string value without touching the
# cell's existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$text) {
    $escaped = $text -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)   # xlPasteValues
}

# Header labels (row 9): "(N=2)" -> "(N=3)"
Set-TextValue $ws.Range("C9") "the % of subjects that this comprises of the total accrual (N=3)"
Set-TextValue $ws.Range("E9") "% of the subjects that this comprises of the total accrual (N=3)"

# Recomputed percentages (count / 3 * 100), formatted like the original table.
Set-TextValue $ws.Range("C10") " 66.67"
Set-TextValue $ws.Range("E10") " 66.67"

Set-TextValue $ws.Range("C11") " 33.33"

Set-TextValue $ws.Range("C12") " 33.33"

Set-TextValue $ws.Range("C13") "100.00"

Set-TextValue $ws.Range("C14") " 66.67"

Set-TextValue $ws.Range("C15") " 66.67"
Set-TextValue $ws.Range("E15") " 33.33"

Set-TextValue $ws.Range("C16") " 33.33"
Set-TextValue $ws.Range("E16") " 33.33"

Set-TextValue $ws.Range("C17") " 66.67"
Set-TextValue $ws.Range("E17") " 66.67"

Set-TextValue $ws.Range("C18") " 66.67"
Set-TextValue $ws.Range("E18") " 33.33"

Set-TextValue $ws.Range("C19") " 33.33"

Set-TextValue $ws.Range("C20") " 33.33"
Set-TextValue $ws.Range("E20") " 33.33"

Set-TextValue $ws.Range("C21") "100.00"
Set-TextValue $ws.Range("E21") "100.00"

Set-TextValue $ws.Range("C22") " 33.33"

Set-TextValue $ws.Range("C23") "100.00"
Set-TextValue $ws.Range("E23") " 33.33"

Set-TextValue $ws.Range("C24") " 66.67"

Set-TextValue $ws.Range("C25") " 33.33"
